# Add the new "Comparison Results" worksheet after Sheet1 by copying Sheet1
# (this preserves the sheet formatting defaults used by the rest of the
# workbook) and then replacing its contents with the new comparison data.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "Comparison Results"
$newSheet.Cells.Clear()

# Header row
$newSheet.Range("A2").Value = "AD Group 1"
$newSheet.Range("B2").Value = "AD Group 2"
$newSheet.Range("C2").Value = "Comparison Results"
$newSheet.Range("D2").Value = "Timestamp"

# Result rows
$newSheet.Range("A3").Value = "GrpGWGCC_LOG_Admins"
$newSheet.Range("B3").Value = "GrpGWGCC_LOG_Users"
$newSheet.Range("C3").Value = "Ramachandran, Hemanathan (g2gyram) (g2gyram); Ramalingam, Karthikeyan (e0hrama) (e0hrama)"
$newSheet.Range("D3").Value = "2024-09-24 13:54:04"

$newSheet.Range("A4").Value = "GrpGWGCC_LOG_Admins"
$newSheet.Range("B4").Value = "GrpGWGCC_LOG_Users"
$newSheet.Range("C4").Value = "Ramachandran, Hemanathan (g2gyram) (g2gyram); Ramalingam, Karthikeyan (e0hrama) (e0hrama)"
$newSheet.Range("D4").Value = "2024-09-24 14:02:15"

# Size the columns to fit their (now much wider) content, as Excel's
# "AutoFit" would when the script populates the sheet.
$newSheet.Columns.Item(1).ColumnWidth = 22.0
$newSheet.Columns.Item(2).ColumnWidth = 20.333333333333332
$newSheet.Columns.Item(3).ColumnWidth = 86.83333333333333
$newSheet.Columns.Item(4).ColumnWidth = 17.333333333333332

[void]$newSheet.Range("F7").Select()
